{"js": "// Update the date line and every \"AxB=C\" multiplication answer in the\n// practice sheet table to the values from the next day's output.\nconst replacements = [\n  [\"2024-12-14 Saturday\", \"2024-12-15 Sunday\"],\n  [\"47\u00d782=3854\", \"15\u00d727=405\"],\n  [\"28\u00d727=756\", \"13\u00d740=520\"],\n  [\"46\u00d747=2162\", \"72\u00d725=1800\"],\n  [\"58\u00d767=3886\", \"31\u00d753=1643\"],\n  [\"98\u00d791=8918\", \"59\u00d781=4779\"],\n  [\"79\u00d756=4424\", \"78\u00d766=5148\"],\n  [\"53\u00d715=795\", \"15\u00d794=1410\"],\n  [\"38\u00d713=494\", \"81\u00d789=7209\"],\n  [\"39\u00d726=1014\", \"54\u00d732=1728\"],\n  [\"11\u00d725=275\", \"21\u00d795=1995\"],\n  [\"76\u00d713=988\", \"95\u00d768=6460\"],\n  [\"71\u00d766=4686\", \"87\u00d717=1479\"],\n  [\"99\u00d760=5940\", \"34\u00d780=2720\"],\n  [\"69\u00d761=4209\", \"98\u00d716=1568\"],\n  [\"29\u00d759=1711\", \"16\u00d752=832\"],\n  [\"70\u00d788=6160\", \"54\u00d731=1674\"],\n  [\"72\u00d740=2880\", \"24\u00d732=768\"],\n  [\"61\u00d746=2806\", \"26\u00d729=754\"],\n  [\"57\u00d750=2850\", \"36\u00d774=2664\"],\n  [\"62\u00d729=1798\", \"54\u00d756=3024\"],\n  [\"89\u00d793=8277\", \"48\u00d717=816\"],\n  [\"69\u00d721=1449\", \"63\u00d765=4095\"],\n  [\"25\u00d737=925\", \"28\u00d712=336\"],\n  [\"40\u00d736=1440\", \"59\u00d777=4543\"],\n  [\"35\u00d746=1610\", \"99\u00d747=4653\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"AxB=C\" multiplication answer in the\n# practice sheet table to the values from the next day's output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-14 Saturday\", \"2024-12-15 Sunday\"),\n    @(\"47\u00d782=3854\", \"15\u00d727=405\"),\n    @(\"28\u00d727=756\", \"13\u00d740=520\"),\n    @(\"46\u00d747=2162\", \"72\u00d725=1800\"),\n    @(\"58\u00d767=3886\", \"31\u00d753=1643\"),\n    @(\"98\u00d791=8918\", \"59\u00d781=4779\"),\n    @(\"79\u00d756=4424\", \"78\u00d766=5148\"),\n    @(\"53\u00d715=795\", \"15\u00d794=1410\"),\n    @(\"38\u00d713=494\", \"81\u00d789=7209\"),\n    @(\"39\u00d726=1014\", \"54\u00d732=1728\"),\n    @(\"11\u00d725=275\", \"21\u00d795=1995\"),\n    @(\"76\u00d713=988\", \"95\u00d768=6460\"),\n    @(\"71\u00d766=4686\", \"87\u00d717=1479\"),\n    @(\"99\u00d760=5940\", \"34\u00d780=2720\"),\n    @(\"69\u00d761=4209\", \"98\u00d716=1568\"),\n    @(\"29\u00d759=1711\", \"16\u00d752=832\"),\n    @(\"70\u00d788=6160\", \"54\u00d731=1674\"),\n    @(\"72\u00d740=2880\", \"24\u00d732=768\"),\n    @(\"61\u00d746=2806\", \"26\u00d729=754\"),\n    @(\"57\u00d750=2850\", \"36\u00d774=2664\"),\n    @(\"62\u00d729=1798\", \"54\u00d756=3024\"),\n    @(\"89\u00d793=8277\", \"48\u00d717=816\"),\n    @(\"69\u00d721=1449\", \"63\u00d765=4095\"),\n    @(\"25\u00d737=925\", \"28\u00d712=336\"),\n    @(\"40\u00d736=1440\", \"59\u00d777=4543\"),\n    @(\"35\u00d746=1610\", \"99\u00d747=4653\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
